{"js": "// Each (oldText -> newText) pair is unique in the document, so a\n// matchCase search for the old text followed by a replace-insert\n// targets exactly the run that needs updating.\nconst replacements = [\n  [\"2024-07-11 Thursday\", \"2024-07-12 Friday\"],\n  [\"22\u00d778=1716\", \"58\u00d728=1624\"],\n  [\"36\u00d747=1692\", \"51\u00d781=4131\"],\n  [\"44\u00d715=660\", \"94\u00d748=4512\"],\n  [\"44\u00d770=3080\", \"69\u00d772=4968\"],\n  [\"29\u00d736=1044\", \"14\u00d761=854\"],\n  [\"27\u00d738=1026\", \"71\u00d794=6674\"],\n  [\"65\u00d783=5395\", \"47\u00d743=2021\"],\n  [\"27\u00d732=864\", \"68\u00d714=952\"],\n  [\"29\u00d726=754\", \"75\u00d774=5550\"],\n  [\"43\u00d794=4042\", \"42\u00d785=3570\"],\n  [\"55\u00d750=2750\", \"71\u00d751=3621\"],\n  [\"29\u00d794=2726\", \"71\u00d784=5964\"],\n  [\"68\u00d796=6528\", \"18\u00d715=270\"],\n  [\"14\u00d789=1246\", \"56\u00d718=1008\"],\n  [\"95\u00d779=7505\", \"98\u00d724=2352\"],\n  [\"31\u00d721=651\", \"85\u00d769=5865\"],\n  [\"85\u00d728=2380\", \"76\u00d713=988\"],\n  [\"67\u00d794=6298\", \"59\u00d744=2596\"],\n  [\"61\u00d717=1037\", \"61\u00d732=1952\"],\n  [\"54\u00d799=5346\", \"83\u00d763=5229\"],\n  [\"21\u00d713=273\", \"92\u00d726=2392\"],\n  [\"12\u00d756=672\", \"42\u00d752=2184\"],\n  [\"90\u00d799=8910\", \"95\u00d722=2090\"],\n  [\"74\u00d720=1480\", \"14\u00d761=854\"],\n  [\"37\u00d712=444\", \"62\u00d771=4402\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each (old -> new) pair below is unique in the document, so a simple\n# MatchCase / whole-text Find & Replace targets exactly the right run.\n$pairs = @(\n    @('2024-07-11 Thursday', '2024-07-12 Friday'),\n    @('22\u00d778=1716', '58\u00d728=1624'),\n    @('36\u00d747=1692', '51\u00d781=4131'),\n    @('44\u00d715=660', '94\u00d748=4512'),\n    @('44\u00d770=3080', '69\u00d772=4968'),\n    @('29\u00d736=1044', '14\u00d761=854'),\n    @('27\u00d738=1026', '71\u00d794=6674'),\n    @('65\u00d783=5395', '47\u00d743=2021'),\n    @('27\u00d732=864', '68\u00d714=952'),\n    @('29\u00d726=754', '75\u00d774=5550'),\n    @('43\u00d794=4042', '42\u00d785=3570'),\n    @('55\u00d750=2750', '71\u00d751=3621'),\n    @('29\u00d794=2726', '71\u00d784=5964'),\n    @('68\u00d796=6528', '18\u00d715=270'),\n    @('14\u00d789=1246', '56\u00d718=1008'),\n    @('95\u00d779=7505', '98\u00d724=2352'),\n    @('31\u00d721=651', '85\u00d769=5865'),\n    @('85\u00d728=2380', '76\u00d713=988'),\n    @('67\u00d794=6298', '59\u00d744=2596'),\n    @('61\u00d717=1037', '61\u00d732=1952'),\n    @('54\u00d799=5346', '83\u00d763=5229'),\n    @('21\u00d713=273', '92\u00d726=2392'),\n    @('12\u00d756=672', '42\u00d752=2184'),\n    @('90\u00d799=8910', '95\u00d722=2090'),\n    @('74\u00d720=1480', '14\u00d761=854'),\n    @('37\u00d712=444', '62\u00d771=4402'),\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
